# adding notices for ensino
# Inserts four new "ensino" i18n rows into the i18n table on Sheet1, expanding
# the table/dimension from A1:G591 to A1:G595, matching the target commit.
#
# New rows (final row numbers after all inserts):
#   562 certified_software_notice_toconline_ensino
#   566 document_certified_notice_toconline_ensino
#   568 document_certified_notice_non_hashed_toconline_ensino
#   570 document_certified_notice_short_toconline_ensino
#
# Strategy: insert a blank row immediately above each target ORIGINAL row
# (567, 566, 565, 562 - processed in that descending order so row numbers used
# in each Insert() call stay valid), fill in the three cell values, then copy
# number-format/style from the most similar neighbouring cells (per-column, not
# per full-row, to avoid the engine synthesizing brand-new blended styles).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# Step 1: new row before original row 567  ->  ends up at row 570
#   document_certified_notice_short_toconline_ensino
# ---------------------------------------------------------------------------
$ws.Rows("567:567").Insert()
$ws.Range("A567").Value = "document_certified_notice_short_toconline_ensino"
$ws.Range("B567").Value = "Emitido por TOConline - https://www.toconline.pt"
$ws.Range("C567").Value = "Processed by TOConline - https://www.toconline.pt"

$ws.Range("A566:C566").Copy()
$ws.Range("A567:C567").PasteSpecial($xlPasteFormats)
$ws.Range("C459").Copy()
$ws.Range("D567:G567").PasteSpecial($xlPasteFormats)
$ws.Rows("567:567").RowHeight = 34

# ---------------------------------------------------------------------------
# Step 2: new row before original row 566  ->  ends up at row 568
#   document_certified_notice_non_hashed_toconline_ensino
# ---------------------------------------------------------------------------
$ws.Rows("566:566").Insert()
$ws.Range("A566").Value = "document_certified_notice_non_hashed_toconline_ensino"
$ws.Range("B566").Value = "Emitido por programa certificado n$([char]0x00BA) 1662/AT - TOConline"
$ws.Range("C566").Value = "Issued by certified program nr. 1662/AT - TOConline"

$ws.Range("A565:C565").Copy()
$ws.Range("A566:C566").PasteSpecial($xlPasteFormats)
$ws.Range("D565:G565").Copy()
$ws.Range("D566:G566").PasteSpecial($xlPasteFormats)
$ws.Rows("566:566").RowHeight = 34

# ---------------------------------------------------------------------------
# Step 3: new row before original row 565  ->  ends up at row 566
#   document_certified_notice_toconline_ensino
# ---------------------------------------------------------------------------
$ws.Rows("565:565").Insert()
$ws.Range("A565").Value = "document_certified_notice_toconline_ensino"
$ws.Range("B565").Value = "Processado por programa certificado n$([char]0x00BA) 1662/AT - TOConline"
$ws.Range("C565").Value = "Processed by certified program nr. 1662/AT - TOConline"

$ws.Range("A564:C564").Copy()
$ws.Range("A565:C565").PasteSpecial($xlPasteFormats)
$ws.Range("D564:G564").Copy()
$ws.Range("D565:G565").PasteSpecial($xlPasteFormats)
$ws.Rows("565:565").RowHeight = 34

# ---------------------------------------------------------------------------
# Step 4: new row before original row 562  ->  ends up at row 562
#   certified_software_notice_toconline_ensino
# ---------------------------------------------------------------------------
$ws.Rows("562:562").Insert()
$ws.Range("A562").Value = "certified_software_notice_toconline_ensino"
$ws.Range("B562").Value = "Emitido por TOConline - https://www.toconline.pt"
$ws.Range("C562").Value = "Created by TOConline - https://www.toconline.pt"

$ws.Range("A561:C561").Copy()
$ws.Range("A562:C562").PasteSpecial($xlPasteFormats)
$ws.Rows("562:562").RowHeight = 17

# ---------------------------------------------------------------------------
# Grow the "i18n" table and sheet dimension to cover the four new rows.
# ---------------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:G595"))

# ---------------------------------------------------------------------------
# Match the saved view state (scroll position + active selection).
# ---------------------------------------------------------------------------
$ws.Range("J595").Select()
$excel.ActiveWindow.ScrollRow = 551
$excel.ActiveWindow.ScrollColumn = 1
